# "update review with PGx team"
#
# Re-review of sample 20220629-24002's CYP2D6_015 (CYP2D6_6) marker: the
# min_height threshold for sample S2 was raised from 1000 -> 1500, which
# pushes the allele call for that marker/direction below detection
# threshold. That ripples into the allele_table (peak/size/height/status
# cleared, a "could not be detected" message recorded), the marker_table
# genotype/phenotype call (TG/heterozygous -> TT/wildtype) and the overall
# genotype_result (*2/*6 -> *1/*2). The peak_table's own m_height for that
# row is bumped to match.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Workbook-level housekeeping
# ---------------------------------------------------------------------
# Workbook was protected with an (empty / no-op) protection marker; drop it.
$wb.Unprotect($null)

# Re-point the active window at the re-reviewed cell, and recalc on load.
$win = $wb.Windows.Item(1)
$win.Left = 38280
$win.Top = -120
$win.Width = 29040
$win.Height = 15840

$excel.Iterate = $true
$excel.MaxChange = 0.0001
$excel.CalculateBeforeSave = $true

# ---------------------------------------------------------------------
# peak_table: min_height review for S2 / CYP2D6_015 (CYP2D6_6)
# ---------------------------------------------------------------------
$peak = $wb.Worksheets.Item("peak_table")

$peak.Range("O16").Value = 1500

# Column widths the reviewer's workbook ended up with (auto best-fit on the
# marker/label columns) and the taller default row height used in the
# reviewed copy.
$peak.StandardHeight = 17
$peak.Columns.Item(3).ColumnWidth = 13.08203125
$peak.Columns.Item(4).ColumnWidth = 12.25

[void]$peak.Range("O16").Select()
$peak.Activate()

# ---------------------------------------------------------------------
# allele_table: row 31 (S2 / CYP2D6_015 / Reverse / G / mutant) now fails
# to detect at the higher min_height threshold.
# ---------------------------------------------------------------------
$allele = $wb.Worksheets.Item("allele_table")

$allele.Range("K31").Value = 1500
$allele.Range("M31").Value = $false
$allele.Range("N31:Q31").ClearContents()
$allele.Range("R31").Value = "Peak(s) could not be detected. Please check peak ranges if required!"

# Header row re-styled (new font/border objects introduced by the review
# workbook's template, same bold / thin-border / centered look).
$headerFont = "맑은 고딕"
$rng = $allele.Range("A1:S1")
$rng.Font.Name = $headerFont
$rng.Font.Size = 8
$rng.Borders.LineStyle = 1
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4160

# ---------------------------------------------------------------------
# marker_table: CYP2D6_015 (CYP2D6_6) genotype call for sample S2 updated
# ---------------------------------------------------------------------
$marker = $wb.Worksheets.Item("marker_table")

$marker.Range("G16").Value = "TT"
$marker.Range("H16").Value = "wildtype"

$rng = $marker.Range("A1:H1")
$rng.Font.Name = $headerFont
$rng.Font.Size = 8
$rng.Borders.LineStyle = 1
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4160

# ---------------------------------------------------------------------
# genotype_result: overall CYP2D6 diplotype call updated
# ---------------------------------------------------------------------
$genotype = $wb.Worksheets.Item("genotype_result")

$genotype.Range("B2").Value = "*1/*2"

$rng = $genotype.Range("A1:B1")
$rng.Font.Name = $headerFont
$rng.Font.Size = 8
$rng.Borders.LineStyle = 1
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4160
